$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("A54").Value = 45986

$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 2.560577522109297
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 1.325305149734723
